# expense ratios method added. stock split amended
# Applies a 9-for-1 stock-split style adjustment (multiply holdings by 9)
# to product_id 3105371 rows whose month_date >= 2024-03-31 (Excel serial 45382).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E205").Value = 129419562.8330666
$ws.Range("E214").Value = 142361519.1163732
$ws.Range("E215").Value = 127349488.3129605
$ws.Range("E216").Value = 132844554.6567569
$ws.Range("E217").Value = 113882959.1049857
$ws.Range("E223").Value = 129419562.8330666
$ws.Range("E224").Value = 127349488.3129605
$ws.Range("E225").Value = 120767776.9606881
$ws.Range("E226").Value = 113882959.1049857
$ws.Range("E227").Value = 101422078.0546259
$ws.Range("E228").Value = 98482205.27843586
$ws.Range("E229").Value = 84996718.23479602
$ws.Range("E232").Value = 129419562.8330666
$ws.Range("E233").Value = 127349488.3129605
$ws.Range("E234").Value = 120767776.9606881
$ws.Range("E235").Value = 113882959.1049857
$ws.Range("E236").Value = 111564285.8600885
$ws.Range("E237").Value = 98482205.27843586
$ws.Range("E238").Value = 84996718.23479602
$ws.Range("E239").Value = 84817029.33548039
$ws.Range("E240").Value = 76516135.42001562
$ws.Range("E241").Value = 76516135.42001562
$ws.Range("E445").Value = 670579715.9906673
$ws.Range("E454").Value = 737637687.589734
$ws.Range("E455").Value = 613196662.7488316
$ws.Range("E456").Value = 606267241.0858063
$ws.Range("E457").Value = 715251887.5900874
$ws.Range("E463").Value = 670579715.9906673
$ws.Range("E464").Value = 613196662.7488316
$ws.Range("E465").Value = 606267241.0858063
$ws.Range("E466").Value = 715251887.5900874
$ws.Range("E467").Value = 755995436.3849871
$ws.Range("E468").Value = 742036702.4947542
$ws.Range("E469").Value = 756691955.8285863
$ws.Range("E472").Value = 670579715.9906673
$ws.Range("E473").Value = 613196662.7488316
$ws.Range("E474").Value = 606267241.0858063
$ws.Range("E475").Value = 786777076.3490964
$ws.Range("E476").Value = 755995436.3849871
$ws.Range("E477").Value = 742036702.4947542
$ws.Range("E478").Value = 756691955.8285863
$ws.Range("E479").Value = 816908391.0171666
$ws.Range("E480").Value = 984761325.1397109
$ws.Range("E481").Value = 984761325.1397109
$ws.Range("E685").Value = 23691528.33574733
$ws.Range("E694").Value = 23691528.33574733
$ws.Range("E695").Value = 21347482.79765468
$ws.Range("E696").Value = 19750402.02934099
$ws.Range("E697").Value = 18375745.21828512
$ws.Range("E703").Value = 23691528.33574733
$ws.Range("E704").Value = 21347482.79765468
$ws.Range("E705").Value = 19750402.02934099
$ws.Range("E706").Value = 18375745.21828512
$ws.Range("E707").Value = 17830641.44426316
$ws.Range("E708").Value = 14595217.40708673
$ws.Range("E709").Value = 14254062.37470174
$ws.Range("E712").Value = 23691528.33574733
$ws.Range("E713").Value = 21347482.79765468
$ws.Range("E714").Value = 19750402.02934099
$ws.Range("E715").Value = 18375745.21828512
$ws.Range("E716").Value = 16209674.04023923
$ws.Range("E717").Value = 14595217.40708673
$ws.Range("E718").Value = 14254062.37470174
$ws.Range("E719").Value = 14204545.44844807
$ws.Range("E720").Value = 12891577.47979644
$ws.Range("E721").Value = 12891577.47979644
$ws.Range("E925").Value = 5321170441.350183
$ws.Range("E934").Value = 5321170441.350183
$ws.Range("E935").Value = 6101795369.01027
$ws.Range("E936").Value = 7144159802.322845
$ws.Range("E937").Value = 7554734637.712518
$ws.Range("E943").Value = 5321170441.350183
$ws.Range("E944").Value = 6101795369.01027
$ws.Range("E945").Value = 7858575782.55513
$ws.Range("E946").Value = 7554734637.712518
$ws.Range("E947").Value = 7821854623.467119
$ws.Range("E948").Value = 7814951383.874064
$ws.Range("E949").Value = 8414861315.225637
$ws.Range("E952").Value = 5321170441.350183
$ws.Range("E953").Value = 6101795369.01027
$ws.Range("E954").Value = 7144159802.322845
$ws.Range("E955").Value = 7554734637.712518
$ws.Range("E956").Value = 7821854623.467119
$ws.Range("E957").Value = 7814951383.874064
$ws.Range("E958").Value = 8414861315.225637
$ws.Range("E959").Value = 8720182095.949892
$ws.Range("E960").Value = 9029484129.394419
$ws.Range("E961").Value = 9029484129.394419
$ws.Range("E1165").Value = 76394307.66843285
$ws.Range("E1174").Value = 76394307.66843285
$ws.Range("E1175").Value = 76408545.2852505
$ws.Range("E1176").Value = 69224361.82157111
$ws.Range("E1177").Value = 61541032.91615139
$ws.Range("E1183").Value = 76394307.66843285
$ws.Range("E1184").Value = 76408545.2852505
$ws.Range("E1185").Value = 69224361.82157111
$ws.Range("E1186").Value = 61541032.91615139
$ws.Range("E1187").Value = 58768272.75972089
$ws.Range("E1188").Value = 51822465.27962469
$ws.Range("E1189").Value = 46661034.85907778
$ws.Range("E1192").Value = 76394307.66843285
$ws.Range("E1193").Value = 76408545.2852505
$ws.Range("E1194").Value = 69224361.82157111
$ws.Range("E1195").Value = 61541032.91615139
$ws.Range("E1196").Value = 64645100.03569299
$ws.Range("E1197").Value = 51822465.27962469
$ws.Range("E1198").Value = 46661034.85907778
$ws.Range("E1199").Value = 41338708.79593902
$ws.Range("E1200").Value = 35595297.54824661
$ws.Range("E1201").Value = 35595297.54824661
